$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last_name on row 7 (Filiberto Christiansen -> Actualizado)
$ws.Range("B7").Value = "Actualizado"

# Delete the row for "Prof. Jovany Daugherty" (row 8) - data has been
# related/imported so this duplicate row is removed and subsequent rows
# shift up.
$ws.Rows(8).Delete()

# Update selection to match final cursor position
$ws.Range("C6").Select()
